$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.672.83"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "3.780.44"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'595.46"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'166.93"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "3.766.86"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "'6.29"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'36.02"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "4.414.62"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "3.757.01"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'18.52"
$ws.Range("D18").Value = "67.605.26"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'10.02"
$ws.Range("E21").Value = "  -5.87%  "
$ws.Range("D22").Value = "'459.66"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "'0.696"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +6.06%  "
$ws.Range("D25").Value = "'83.38"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'10.01"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "'2.23"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").Value = "'7.21"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'29.56"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'3.37"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'0.995"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "'5.77"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D43").Value = "'45.46"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").Value = "'48.07"
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").Value = "'0.299"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "'149.74"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").Value = "'8.31"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "'393.94"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "'26.71"
$ws.Range("E49").Value = "  +6.53%  "
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "2.719.87"
$ws.Range("E51").Value = "  -0.82%  "
